$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$districts = @(
    @(44181, 'Bánovce nad Bebravou', 14),
    @(44181, 'Banská Bystrica', 88),
    @(44181, 'Banská Štiavnica', 4),
    @(44181, 'Bardejov', 63),
    @(44181, 'Bratislava', 213),
    @(44181, 'Brezno', 26),
    @(44181, 'Bytča', 16),
    @(44181, 'Čadca', 45),
    @(44181, 'Detva', 15),
    @(44181, 'Dolný Kubín', 18),
    @(44181, 'Dunajská Streda', 38),
    @(44181, 'Galanta', 29),
    @(44181, 'Gelnica', 13),
    @(44181, 'Hlohovec', 27),
    @(44181, 'Humenné', 12),
    @(44181, 'Ilava', 41),
    @(44181, 'Kežmarok', 16),
    @(44181, 'Komárno', 12),
    @(44181, 'Košice', 144),
    @(44181, 'Košice - okolie', 53),
    @(44181, 'Krupina', 13),
    @(44181, 'Kysucké Nové Mesto', 5),
    @(44181, 'Levice', 31),
    @(44181, 'Levoča', 47),
    @(44181, 'Liptovský Mikuláš', 37),
    @(44181, 'Lučenec', 47),
    @(44181, 'Malacky', 30),
    @(44181, 'Martin', 63),
    @(44181, 'Medzilaborce', 2),
    @(44181, 'Michalovce', 43),
    @(44181, 'Myjava', 59),
    @(44181, 'Námestovo', 7),
    @(44181, 'Nitra', 201),
    @(44181, 'Nové Mesto nad Váhom', 54),
    @(44181, 'Nové Zámky', 54),
    @(44181, 'Partizánske', 3),
    @(44181, 'Pezinok', 13),
    @(44181, 'Piešťany', 34),
    @(44181, 'Poltár', 6),
    @(44181, 'Poprad', 65),
    @(44181, 'Považská Bystrica', 64),
    @(44181, 'Prešov', 97),
    @(44181, 'Prievidza', 64),
    @(44181, 'Púchov', 66),
    @(44181, 'Revúca', 12),
    @(44181, 'Rimavská Sobota', 45),
    @(44181, 'Rožňava', 29),
    @(44181, 'Ružomberok', 26),
    @(44181, 'Sabinov', 13),
    @(44181, 'Senec', 26),
    @(44181, 'Senica', 84),
    @(44181, 'Skalica', 20),
    @(44181, 'Snina', 6),
    @(44181, 'Sobrance', 12),
    @(44181, 'Spišská Nová Ves', 91),
    @(44181, 'Stará Ľubovňa', 34),
    @(44181, 'Stropkov', 12),
    @(44181, 'Svidník', 41),
    @(44181, 'Šaľa', 16),
    @(44181, 'Topoľčany', 25),
    @(44181, 'Trebišov', 70),
    @(44181, 'Trenčín', 96),
    @(44181, 'Trnava', 80),
    @(44181, 'Turčianske Teplice', 2),
    @(44181, 'Tvrdošín', 32),
    @(44181, 'Veľký Krtíš', 44),
    @(44181, 'Vranov nad Topľou', 58),
    @(44181, 'Zlaté Moravce', 18),
    @(44181, 'Zvolen', 45),
    @(44181, 'Žarnovica', 11),
    @(44181, 'Žiar nad Hronom', 23),
    @(44181, 'Žilina', 82)
)

$startRow = 6222
for ($i = 0; $i -lt $districts.Count; $i++) {
    $row = $startRow + $i
    $rec = $districts[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
}

Write-Host "Done writing $($districts.Count) rows"